$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计"
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# IMPORTANT: fetch "总计" only AFTER the insertion above, since inserting a
# sheet shifts everybody's position and a reference captured beforehand
# would otherwise keep pointing at the (now wrong) positional slot.
$totalSheet = $wb.Worksheets.Item("总计")

# Header labels
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($j = 0; $j -lt $headers.Count; $j++) {
    $newSheet.Cells.Item(1, $j + 2).Value = $headers[$j]
}

# Data rows: index, fund code, fund name, fund size, total stock position,
# position ratio, market value held (100M yuan), position rank
$data = @(
    @(0, "011271", "汇添富价值成长均衡投资混合A", "24.78", "92.58", "3.85", "0.9540", 8),
    @(1, "011410", "中信建投量化进取6个月持有期混合A", "9.13", "93.80", "0.79", "0.0721", 8),
    @(2, "009954", "北信瑞丰优选成长股票", "0.57", "94.37", "4.96", "0.0283", 7),
    @(3, "001829", "北信瑞丰中国智造主题灵活配置混合", "0.36", "94.06", "5.07", "0.0183", 6),
    @(4, "011411", "中信建投量化进取6个月持有期混合C", "2.15", "93.80", "0.79", "0.0170", 8),
    @(5, "011272", "汇添富价值成长均衡投资混合C", "0.32", "92.58", "3.85", "0.0123", 8),
    @(6, "002123", "北信瑞丰外延增长主题灵活配置混合", "0.17", "94.48", "5.32", "0.0090", 6),
    @(7, "004726", "先锋聚优灵活配置混合A", "0.06", "93.23", "5.02", "0.0030", 9),
    @(8, "004727", "先锋聚优灵活配置混合C", "0.04", "93.23", "5.02", "0.0020", 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $newSheet.Cells.Item($r, 1).Value = $row[0]

    # Fund code must stay textual (keep leading zeros)
    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    # Numeric-looking columns that must be preserved verbatim as text
    # (keep trailing zeros such as "93.80" / "0.9540")
    $sizeCell = $newSheet.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[3]

    $posCell = $newSheet.Cells.Item($r, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[4]

    $ratioCell = $newSheet.Cells.Item($r, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $row[5]

    $mvCell = $newSheet.Cells.Item($r, 7)
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# Match formatting of the "2021-Q4" sheet: bold/centered/bordered header row
# and bold/centered/bordered index column (A).
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q4Sheet.Range("A2:A10").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)



# ---------------------------------------------------------------------------
# 2. Update "总计" sheet: add a new row for 2022-Q1 above the 2021-Q4 row
# ---------------------------------------------------------------------------
$oldIndex = $totalSheet.Cells.Item(2, 1).Value2
$oldDate = $totalSheet.Cells.Item(2, 2).Value2
$oldCount = $totalSheet.Cells.Item(2, 3).Value2
$oldValue = $totalSheet.Cells.Item(2, 4).Value2

# Push the existing 2021-Q4 row down to row 3
$totalSheet.Cells.Item(3, 1).Value = $oldIndex + 1
$totalSheet.Cells.Item(3, 2).Value = $oldDate
$totalSheet.Cells.Item(3, 3).Value = $oldCount
$totalSheet.Cells.Item(3, 4).Value = $oldValue

# Copy the index-column formatting down to the row that now holds it
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# Write the new 2022-Q1 row into row 2
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 1.12

Write-Host "done"
